$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 302; this shifts every row from the old
# 302 onward down by one (old 302 -> 303, ..., old 361 -> 362), and
# inherits the formatting (e.g. the date style on column D) from the
# row above, exactly like a manual Excel "Insert Row" would.
$ws.Rows.Item(302).Insert()

# Populate the new row 302 with the new weekly data point. All
# non-changed fields mirror what used to be in row 302 (now row 303);
# only Fecha (D), Volumen (J) and Origen (O) differ.
$ws.Cells.Item(302, 1).Value = 10
$ws.Cells.Item(302, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(302, 3).Value = "La Araucanía"
$ws.Cells.Item(302, 4).Value = "2023-03-10"
$ws.Cells.Item(302, 5).Value = 9
$ws.Cells.Item(302, 6).Value = 100112052
$ws.Cells.Item(302, 7).Value = "Albahaca"
$ws.Cells.Item(302, 8).Value = "Sin especificar"
$ws.Cells.Item(302, 9).Value = "Primera"
$ws.Cells.Item(302, 10).Value = 35
$ws.Cells.Item(302, 11).Value = 6000
$ws.Cells.Item(302, 12).Value = 6000
$ws.Cells.Item(302, 13).Value = 6000
$ws.Cells.Item(302, 14).Value = "`$/paquete"
$ws.Cells.Item(302, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(302, 16).Value = 6000
$ws.Cells.Item(302, 17).Value = 1
$ws.Cells.Item(302, 18).Value = "Hortaliza"
